$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.122.76'
$ws.Range("E2").Value = '  -4.56%  '

$ws.Range("D3").Value = '1.651.59'
$ws.Range("E3").Value = '  -3.66%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.06'
$ws.Range("D5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5113'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.00%  '

$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2597'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06443'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.96%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.78'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.47%  '

$ws.Range("E11").Value = '  -0.36%  '

$ws.Range("D12").Value = '1.669.64'
$ws.Range("E12").Value = '  -2.47%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.282'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -4.01%  '

$ws.Range("D14").Value = '1.878.92'
$ws.Range("E14").Value = '  -3.66%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5497'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.85%  '

$ws.Range("D16").Value = '0.0₅8020'
$ws.Range("E16").Value = '  -2.04%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.93'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -5.61%  '

$ws.Range("D18").Value = '26.135.35'
$ws.Range("E18").Value = '  -4.55%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.010'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.32'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -5.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.400'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.07'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.39%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.034'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.010'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.863'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +8.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '143.57'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1171'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.74%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.933'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.85%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.86'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05105'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.92%  '

$ws.Range("E31").Value = '  -3.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.352'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.49%  '

$ws.Range("E33").Value = '  -3.43%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.553'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -5.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.351'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9169'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.653'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5718'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.61%  '

$ws.Range("D39").Value = '1.154.71'
$ws.Range("E39").Value = '  -2.67%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01578'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.31%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.571'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.010'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.670'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.35%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8282'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '100.30'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.88%  '

$ws.Range("D46").Value = '1.789.88'
$ws.Range("E46").Value = '  -3.68%  '

$ws.Range("D47").Value = '0.0₈114'
$ws.Range("E47").Value = '  -3.24%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4547'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.36%  '

$ws.Range("E49").Value = '  +0.62%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '55.25'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.88%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.820'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.62%  '
